$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 ("5_hypo_dl") was removed entirely from the results table.
$ws.Rows(6).Delete()

# Updated cosinor statistics for the remaining rows (2-5) to match the paper re-run.
$ws.Range("E2").Value = 24.04000000000032
$ws.Range("G2").Value = [double]"1.110223024625157e-16"
$ws.Range("H2").Value = [double]"2.328121676802425e-16"
$ws.Range("I2").Value = ""
$ws.Range("K2").Value = 1.022353704386827
$ws.Range("L2").Value = "[0.9181954150184453, 1.1265119937552086]"
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 1.465447624197041
$ws.Range("P2").Value = "[1.3522370781217328, 1.578658170272349]"
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 36.80528711476322
$ws.Range("T2").Value = "[36.73821276581787, 36.872361463708565]"
$ws.Range("W2").Value = 18.43307307307332
$ws.Range("X2").Value = 17.99991991992016
$ws.Range("Y2").Value = 18.86622622622648
$ws.Range("C3").Value = "2_induction_dd"
$ws.Range("E3").Value = 23.84000000000029
$ws.Range("H3").Value = [double]"2.328121676802425e-16"
$ws.Range("I3").Value = ""
$ws.Range("K3").Value = 1.075783570129403
$ws.Range("L3").Value = "[0.9700216850831822, 1.1815454551756233]"
$ws.Range("O3").Value = 2.283079345852042
$ws.Range("P3").Value = "[2.1824477493406578, 2.383710942363426]"
$ws.Range("S3").Value = 36.41837031685424
$ws.Range("T3").Value = "[36.36465073637344, 36.472089897335046]"
$ws.Range("W3").Value = 15.1774174174176
$ws.Range("X3").Value = 14.79559559559577
$ws.Range("Y3").Value = 15.55923923923942
$ws.Range("C4").Value = "3_hypo_dd"
$ws.Range("E4").Value = 23.70000000000027
$ws.Range("H4").Value = [double]"2.328121676802425e-16"
$ws.Range("I4").Value = [double]"1.110223024625157e-16"
$ws.Range("K4").Value = 0.6375504901043954
$ws.Range("L4").Value = "[0.526386709649568, 0.7487142705592227]"
$ws.Range("O4").Value = 1.427710775505271
$ws.Range("P4").Value = "[1.2390265320464247, 1.6163950189641172]"
$ws.Range("S4").Value = 36.3021525139803
$ws.Range("T4").Value = "[36.23209958631354, 36.372205441647054]"
$ws.Range("W4").Value = 18.31471471471492
$ws.Range("X4").Value = 17.6030030030032
$ws.Range("Y4").Value = 19.02642642642664
$ws.Range("C5").Value = "4_hypo_dl"
$ws.Range("E5").Value = 25.53000000000055
$ws.Range("H5").Value = [double]"2.328121676802425e-16"
$ws.Range("I5").Value = ""
$ws.Range("K5").Value = 1.073489080837021
$ws.Range("L5").Value = "[0.9532821238608733, 1.1936960378131687]"
$ws.Range("O5").Value = -1.899421384152387
$ws.Range("P5").Value = "[-2.0126319302276947, -1.7862108380770803]"
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 36.26360265019388
$ws.Range("T5").Value = "[36.195665564661525, 36.33153973572624]"
$ws.Range("W5").Value = 7.717777777777947
$ws.Range("X5").Value = 7.25777777777794
$ws.Range("Y5").Value = 8.177777777777955
